$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.117.12'
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("D3").Value = '1.789.03'
$ws.Range("E3").Value = '  -2.60%  '
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("E6").Value = '  -1.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.56'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("E9").Value = '  -2.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0704'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.01%  '
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").Value = '2.049.71'
$ws.Range("E12").Value = '  -2.27%  '
$ws.Range("D13").Value = '1.786.82'
$ws.Range("E13").Value = '  -2.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.624'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.33%  '
$ws.Range("D16").Value = '34.101.68'
$ws.Range("E16").Value = '  -1.54%  '
$ws.Range("E17").Value = '  -4.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.96'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.98'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.34%  '
$ws.Range("D20").Value = '0.0₃0783'
$ws.Range("E20").Value = '  -2.45%  '
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.71'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.77%  '
$ws.Range("E23").Value = '  -4.72%  '
$ws.Range("E24").Value = '  -2.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.28'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.40%  '
$ws.Range("E28").Value = '  -2.11%  '
$ws.Range("E29").Value = '  +0.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0514'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.04%  '
$ws.Range("E31").Value = '  +0.57%  '
$ws.Range("E32").Value = '  -4.06%  '
$ws.Range("E33").Value = '  -3.61%  '
$ws.Range("E34").Value = '  -6.30%  '
$ws.Range("D35").Value = '1.394.41'
$ws.Range("E35").Value = '  -3.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.644'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.81%  '
$ws.Range("E37").Value = '  -2.34%  '
$ws.Range("E38").Value = '  -3.75%  '
$ws.Range("E39").Value = '  +1.67%  '
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.49%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.70'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.78%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.913'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.47%  '
$ws.Range("B44").Value = 'BabyDogeCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D44").Value = '0.0₆0147'
$ws.Range("E44").Value = '  +16.38%  '
$ws.Range("E45").Value = '  +1.19%  '
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '107.64'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.89'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.49%  '
$ws.Range("D49").Value = '1.947.43'
$ws.Range("E49").Value = '  -1.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '12.19'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.04%  '
